$wb = $excel.ActiveWorkbook


# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 342
$ws.Range("I5").Value = 109.5
$ws.Range("J5").Value = 574.5
$ws.Range("K5").Value = 109.5
$ws.Range("L5").Value = 574.5
$ws.Range("M5").Value = 5.5
$ws.Range("N5").Value = -804.5
$ws.Range("H9").Value = 656.5333000000001
$ws.Range("I9").Value = 718.9167
$ws.Range("J9").Value = 407
$ws.Range("K9").Value = 718.9167
$ws.Range("L9").Value = 407
$ws.Range("M9").Value = -549.9167
$ws.Range("N9").Value = -745
$ws.Range("H10").Value = 2004
$ws.Range("I10").Value = 2004
$ws.Range("K10").Value = 2004
$ws.Range("M10").Value = -1711
$ws.Range("H13").Value = 1203.75
$ws.Range("I13").Value = 299.5
$ws.Range("J13").Value = 2108
$ws.Range("K13").Value = 299.5
$ws.Range("L13").Value = 2108
$ws.Range("M13").Value = -130.5
$ws.Range("N13").Value = -2446
$ws.Range("H21").Value = 3800
$ws.Range("I21").Value = 3250
$ws.Range("K21").Value = 3250
$ws.Range("M21").Value = -2782
$ws.Range("H23").Value = 3800
$ws.Range("I23").Value = 3250
$ws.Range("K23").Value = 3250
$ws.Range("M23").Value = -3016
$ws.Range("H33").Value = 418.8
$ws.Range("I33").Value = 391.07693
$ws.Range("K33").Value = 391.07693
$ws.Range("M33").Value = -162.07693
$ws.Range("H47").Value = 49666.668
$ws.Range("I47").Value = 25000
$ws.Range("K47").Value = 25000
$ws.Range("M47").Value = -24028
$ws.Range("H55").Value = 545.8570999999999
$ws.Range("I55").Value = 269.5
$ws.Range("J55").Value = 914.3333
$ws.Range("K55").Value = 269.5
$ws.Range("L55").Value = 914.3333
$ws.Range("M55").Value = -55.5
$ws.Range("N55").Value = -1342.3333
$ws.Range("H76").Value = 3524
$ws.Range("I76").Value = 3899
$ws.Range("J76").Value = 3399
$ws.Range("K76").Value = 3899
$ws.Range("L76").Value = 3399
$ws.Range("M76").Value = -3584
$ws.Range("N76").Value = -4029
$ws.Range("H79").Value = 3524
$ws.Range("I79").Value = 3899
$ws.Range("J79").Value = 3399
$ws.Range("K79").Value = 3899
$ws.Range("L79").Value = 3399
$ws.Range("M79").Value = -2807
$ws.Range("N79").Value = -5583
$ws.Range("H80").Value = 1715.4667
$ws.Range("I80").Value = 2032.8334
$ws.Range("J80").Value = 1503.8889
$ws.Range("K80").Value = 6098.5002
$ws.Range("L80").Value = 4511.6667
$ws.Range("M80").Value = -5100.5002
$ws.Range("N80").Value = -6507.6667
$ws.Range("H83").Value = 1715.4667
$ws.Range("I83").Value = 2032.8334
$ws.Range("J83").Value = 1503.8889
$ws.Range("K83").Value = 18295.5006
$ws.Range("L83").Value = 13535.0001
$ws.Range("M83").Value = -13303.5006
$ws.Range("N83").Value = -23519.0001
$ws.Range("H96").Value = 1677.1538
$ws.Range("I96").Value = 2141.1667
$ws.Range("J96").Value = 1279.4286
$ws.Range("K96").Value = 6423.500100000001
$ws.Range("L96").Value = 3838.2858
$ws.Range("M96").Value = -5050.500100000001
$ws.Range("N96").Value = -6584.2858
$ws.Range("H98").Value = 1089.7368
$ws.Range("I98").Value = 408.06668
$ws.Range("K98").Value = 408.06668
$ws.Range("M98").Value = 1089.93332
$ws.Range("H107").Value = 1201.9445
$ws.Range("I107").Value = 1168.7693
$ws.Range("K107").Value = 1168.7693
$ws.Range("M107").Value = 751.2307000000001
$ws.Range("H122").Value = 1089.7368
$ws.Range("I122").Value = 408.06668
$ws.Range("K122").Value = 1224.20004
$ws.Range("M122").Value = 1225.79996
$ws.Range("H123").Value = 84780
$ws.Range("J123").Value = 84780
$ws.Range("L123").Value = 84780
$ws.Range("N123").Value = -94580
$ws.Range("H132").Value = 7382
$ws.Range("I132").Value = 7580.5
$ws.Range("K132").Value = 22741.5
$ws.Range("M132").Value = -20211.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1528.4242
$ws.Range("H131").Value = 39997.5
$ws.Range("J131").Value = 39997.5
$ws.Range("L131").Value = 39997.5
$ws.Range("N131").Value = -50077.5
$ws.Range("H132").Value = 3558.0625
$ws.Range("I132").Value = 1992.2222
$ws.Range("J132").Value = 5571.2856
$ws.Range("K132").Value = 5976.6666
$ws.Range("L132").Value = 16713.8568
$ws.Range("M132").Value = -3446.6666
$ws.Range("N132").Value = -21773.8568
$ws.Range("H135").Value = 61249.25
$ws.Range("J135").Value = 61249.25
$ws.Range("L135").Value = 61249.25
$ws.Range("N135").Value = -71389.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1826.1
$ws.Range("I20").Value = 1826.1
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 1826.1
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -1579.1
$ws.Range("N20").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 5603.645
$ws.Range("I22").Value = 1478.0667
$ws.Range("J22").Value = 9471.375
$ws.Range("K22").Value = 1478.0667
$ws.Range("L22").Value = 9471.375
$ws.Range("M22").Value = -1128.0667
$ws.Range("N22").Value = -10171.375
$ws.Range("H41").Value = 15802
$ws.Range("I41").Value = 15802
$ws.Range("K41").Value = 15802
$ws.Range("M41").Value = -15374
$ws.Range("H59").Value = 24125
$ws.Range("I59").Value = 27000
$ws.Range("J59").Value = 15500
$ws.Range("K59").Value = 27000
$ws.Range("L59").Value = 15500
$ws.Range("M59").Value = -25855
$ws.Range("N59").Value = -17790
$ws.Range("H62").Value = 2836.1428
$ws.Range("I62").Value = 2939
$ws.Range("J62").Value = 2795
$ws.Range("K62").Value = 2939
$ws.Range("L62").Value = 2795
$ws.Range("M62").Value = -2315
$ws.Range("N62").Value = -4043
$ws.Range("H65").Value = 2836.1428
$ws.Range("I65").Value = 2939
$ws.Range("J65").Value = 2795
$ws.Range("K65").Value = 14695
$ws.Range("L65").Value = 13975
$ws.Range("M65").Value = -11575
$ws.Range("N65").Value = -20215
$ws.Range("H93").Value = 20407
$ws.Range("I93").Value = 20407
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 20407
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -18535
$ws.Range("N93").ClearContents()
$ws.Range("H134").Value = 3542.7646
$ws.Range("I134").Value = 3388.3125
$ws.Range("K134").Value = 10164.9375
$ws.Range("M134").Value = -7629.9375

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 448.5
$ws.Range("I122").Value = 448.5
$ws.Range("K122").Value = 4036.5
$ws.Range("M122").Value = -1586.5
$ws.Range("H140").Value = 835048.5
$ws.Range("I140").Value = 835048.5
$ws.Range("K140").Value = 2505145.5
$ws.Range("M140").Value = -2499965.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H70").Value = 7245.0586
$ws.Range("I70").Value = 5996
$ws.Range("K70").Value = 5996
$ws.Range("M70").Value = -5726
$ws.Range("H73").Value = 7245.0586
$ws.Range("I73").Value = 5996
$ws.Range("K73").Value = 5996
$ws.Range("M73").Value = -5060
$ws.Range("H128").Value = 40495
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 40495
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 40495
$ws.Range("M128").ClearContents()
$ws.Range("N128").Value = -50455

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1626.8572
$ws.Range("I16").Value = 1723
$ws.Range("K16").Value = 1723
$ws.Range("M16").Value = -1553
$ws.Range("H29").Value = 2833.3333
$ws.Range("I29").Value = 1750
$ws.Range("K29").Value = 1750
$ws.Range("M29").Value = -1455
$ws.Range("H61").Value = 2089.5881
$ws.Range("I61").Value = 1737.6154
$ws.Range("J61").Value = 3233.5
$ws.Range("K61").Value = 1737.6154
$ws.Range("L61").Value = 3233.5
$ws.Range("M61").Value = -1535.6154
$ws.Range("N61").Value = -3637.5
$ws.Range("H68").Value = 7170.875
$ws.Range("I68").Value = 5142.25
$ws.Range("K68").Value = 5142.25
$ws.Range("M68").Value = -4393.25
$ws.Range("H71").Value = 7170.875
$ws.Range("I71").Value = 5142.25
$ws.Range("K71").Value = 25711.25
$ws.Range("M71").Value = -21967.25
$ws.Range("H93").Value = 2133.5
$ws.Range("I93").Value = 1475.25
$ws.Range("K93").Value = 1475.25
$ws.Range("M93").Value = -227.25
$ws.Range("H113").Value = 2089.5881
$ws.Range("I113").Value = 1737.6154
$ws.Range("J113").Value = 3233.5
$ws.Range("K113").Value = 1737.6154
$ws.Range("L113").Value = 3233.5
$ws.Range("M113").Value = 432.3846000000001
$ws.Range("N113").Value = -7573.5
$ws.Range("H132").Value = 2646.9666
$ws.Range("I132").Value = 2585.5217
$ws.Range("K132").Value = 7756.5651
$ws.Range("M132").Value = -5226.5651

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 4005379.8
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H69").Value = 30343.334
$ws.Range("J69").Value = 30343.334
$ws.Range("L69").Value = 30343.334
$ws.Range("N69").Value = -31841.334
$ws.Range("H72").Value = 30343.334
$ws.Range("J72").Value = 30343.334
$ws.Range("L72").Value = 91030.00199999999
$ws.Range("N72").Value = -98518.00199999999
$ws.Range("H133").Value = 54997
$ws.Range("I133").Value = 49994
$ws.Range("K133").Value = 49994
$ws.Range("M133").Value = -44934
